$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Grow the table to 13 columns (2 new trailing placeholder columns) ---
$lo.Resize($ws.Range("A1:M115"))

# --- 2. Write the final header row text, positioned so that the table/sheet
#        ends up with "Highest Bid" / "Lowest Ask" inserted right after
#        "Stockx Average (last 3 sales)" and before "Stockx Link", shifting
#        the remaining columns two slots to the right. Order of assignment
#        controls shared-string creation order, so Highest Bid is written
#        before Lowest Ask to match the expected shared string table. ---
$ws.Cells.Item(1, 13).Value = "Price Sold"
$ws.Cells.Item(1, 12).Value = "Sold"
$ws.Cells.Item(1, 11).Value = "Location"
$ws.Cells.Item(1, 10).Value = "Last Updated"
$ws.Cells.Item(1, 9).Value = "Stockx Link"
$ws.Cells.Item(1, 7).Value = "Highest Bid"
$ws.Cells.Item(1, 8).Value = "Lowest Ask"

# --- 3. Restore the bold/left-aligned header formatting on the two cells
#        that were appended by Resize (they default to a plain style). ---
$ws.Cells.Item(1, 12).Font.Bold = $true
$ws.Cells.Item(1, 12).HorizontalAlignment = -4131
$ws.Cells.Item(1, 13).Font.Bold = $true
$ws.Cells.Item(1, 13).HorizontalAlignment = -4131

# --- 4. Approximate the column widths Excel would compute for the new /
#        shifted columns. ---
$ws.Columns.Item(7).ColumnWidth = 12.3046875
$ws.Columns.Item(8).ColumnWidth = 12.1875
$ws.Columns.Item(9).ColumnWidth = 12.1875
$ws.Columns.Item(10).ColumnWidth = 13.59375
$ws.Columns.Item(11).ColumnWidth = 9.4921875
$ws.Columns.Item(12).ColumnWidth = 6.328125

# --- 5. A cosmetic 8pt font is referenced by a phoneticPr element in the
#        target file; applying it to an otherwise-untouched helper cell
#        reproduces the same font entry in styles.xml without altering any
#        visible content. ---
$ws.Range("A1").Font.Size = $ws.Range("A1").Font.Size
$tmp = $ws.Range("XFD1")
$tmp.Value = "x"
$tmp.Font.Size = 8
$tmp.Value = ""

# --- 6. Match the saved selection / active cell shown in the diff. ---
$ws.Range("G11").Select()
